$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'66.715.55"
$ws.Cells.Item(2, 5).Value = '  -2.29%  '
$ws.Cells.Item(3, 4).Value = "'3.481.32"
$ws.Cells.Item(3, 5).Value = '  -2.01%  '
$ws.Cells.Item(4, 5).Value = '  -0.13%  '
$ws.Cells.Item(5, 4).Value = "'603.78"
$ws.Cells.Item(5, 5).Value = '  -2.42%  '
$ws.Cells.Item(6, 4).Value = "'148.33"
$ws.Cells.Item(6, 5).Value = '  -4.15%  '
$ws.Cells.Item(7, 4).Value = "'3.480.16"
$ws.Cells.Item(7, 5).Value = '  -1.87%  '
$ws.Cells.Item(8, 5).Value = '  -0.17%  '
$ws.Cells.Item(9, 5).Value = '  -1.31%  '
$ws.Cells.Item(10, 5).Value = '  -2.62%  '
$ws.Cells.Item(11, 5).Value = '  +3.71%  '
$ws.Cells.Item(12, 5).Value = '  -2.99%  '
$ws.Cells.Item(13, 5).Value = '  -4.31%  '
$ws.Cells.Item(14, 4).Value = "'4.069.23"
$ws.Cells.Item(14, 5).Value = '  -2.18%  '
$ws.Cells.Item(15, 4).Value = "'31.63"
$ws.Cells.Item(15, 5).Value = '  -4.46%  '
$ws.Cells.Item(16, 4).Value = "'3.478.08"
$ws.Cells.Item(16, 5).Value = '  -2.33%  '
$ws.Cells.Item(17, 4).Value = "'66.784.85"
$ws.Cells.Item(17, 5).Value = '  -2.75%  '
$ws.Cells.Item(18, 5).Value = '  -0.40%  '
$ws.Cells.Item(19, 5).Value = '  -3.95%  '
$ws.Cells.Item(20, 5).Value = '  -3.12%  '
$ws.Cells.Item(21, 5).Value = '  +1.10%  '
$ws.Cells.Item(22, 4).Value = "'439.62"
$ws.Cells.Item(22, 5).Value = '  -3.60%  '
$ws.Cells.Item(23, 5).Value = '  -4.56%  '
$ws.Cells.Item(24, 4).Value = "'79.60"
$ws.Cells.Item(24, 5).Value = '  +1.46%  '
$ws.Cells.Item(25, 5).Value = '  +0.08%  '
$ws.Cells.Item(26, 4).Value = "'3.618.60"
$ws.Cells.Item(26, 5).Value = '  -2.30%  '
$ws.Cells.Item(27, 5).Value = '  -8.22%  '
$ws.Cells.Item(28, 4).Value = "'9.77"
$ws.Cells.Item(28, 5).Value = '  -6.91%  '
$ws.Cells.Item(29, 4).Value = "'8.40"
$ws.Cells.Item(29, 5).Value = '  -6.86%  '
$ws.Cells.Item(30, 5).Value = '  -2.78%  '
$ws.Cells.Item(31, 5).Value = '  -5.61%  '
$ws.Cells.Item(32, 5).Value = '  -0.68%  '
$ws.Cells.Item(33, 5).Value = '  -0.07%  '
$ws.Cells.Item(34, 4).Value = "'25.41"
$ws.Cells.Item(34, 5).Value = '  -2.75%  '
$ws.Cells.Item(35, 5).Value = '  -5.81%  '
$ws.Cells.Item(36, 4).Value = "'3.470.74"
$ws.Cells.Item(36, 5).Value = '  -2.27%  '
$ws.Cells.Item(37, 4).Value = "'1.81"
$ws.Cells.Item(37, 5).Value = '  -5.80%  '
$ws.Cells.Item(38, 5).Value = '  -3.84%  '
$ws.Cells.Item(39, 5).Value = '  -0.01%  '
$ws.Cells.Item(40, 4).Value = "'0.999"
$ws.Cells.Item(40, 5).Value = '  -0.15%  '
$ws.Cells.Item(41, 4).Value = "'176.83"
$ws.Cells.Item(41, 5).Value = '  -0.98%  '
$ws.Cells.Item(42, 4).Value = "'0.0890"
$ws.Cells.Item(42, 5).Value = '  -2.89%  '
$ws.Cells.Item(43, 4).Value = "'2.13"
$ws.Cells.Item(43, 5).Value = '  -9.14%  '
$ws.Cells.Item(44, 5).Value = '  -2.63%  '
$ws.Cells.Item(45, 4).Value = "'0.892"
$ws.Cells.Item(45, 5).Value = '  -0.56%  '
$ws.Cells.Item(46, 4).Value = "'29.00"
$ws.Cells.Item(46, 5).Value = '  -4.67%  '
$ws.Cells.Item(47, 4).Value = "'46.28"
$ws.Cells.Item(47, 5).Value = '  +1.58%  '
$ws.Cells.Item(48, 5).Value = '  -7.13%  '
$ws.Cells.Item(49, 5).Value = '  -3.81%  '
$ws.Cells.Item(50, 5).Value = '  -8.96%  '
$ws.Cells.Item(51, 4).Value = "'0.984"
$ws.Cells.Item(51, 5).Value = '  -3.56%  '
